$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for rows 2 and 3
$values = @{
    "D" = 0.0161
    "E" = -0.0348
    "G" = 0.09230769230769231
    "H" = 0.09230769230769231
    "I" = 0.08951048951048951
    "J" = 0.08951048951048951
    "K" = 0.109
    "L" = 0.07622377622377623
    "M" = 0.039668
    "N" = 0.01888952380952381
    "O" = 0.3639266055045872
    "P" = 0.039668
    "Q" = 0.01888952380952381
    "R" = 0.3639266055045872
    "U" = 0.274
    "V" = 0.1304761904761905
    "W" = 0.06488095238095239
    "X" = 0.05090557174178309
    "Y" = 0.0139753806391693
    "Z" = 0.9558823529411764
    "AA" = 0.08556149732620322
    "AB" = 0.05090557174178309
    "AC" = 0.03465592558442013
    "AG" = -0.274
    "AJ" = -0.1500547645125958
    "AK" = -0.1664641555285541
    "AP" = -1.412371134020619
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range($col + "2").Value = $val
    $ws.Range($col + "3").Value = $val
}
